$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("A7").Value = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("D7").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E7").Value = "PLACA_20250717165933"
$ws.Range("F7").Value = "Beja"
$ws.Range("G7").Value = 38.02035
$ws.Range("H7").Value = -7.94715
$ws.Range("I7").Value = "962,713,1006,765"
$ws.Range("J7").Value = "'0.76"

# Row 8
$ws.Range("A8").Value = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
$ws.Range("I8").Value = "967,614,1002,659"
$ws.Range("J8").Value = "'0.73"

# Row 9
$ws.Range("A9").Value = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
$ws.Range("I9").Value = "'702,633,740,690"
$ws.Range("J9").Value = "'0.72"

# Row 10
$ws.Range("A10").Value = "dfd476d4-7689-4671-a076-78fe3ce806bb"
$ws.Range("I10").Value = "1254,850,1294,895"
$ws.Range("J10").Value = "'0.67"

# Row 11
$ws.Range("A11").Value = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D11").Value = "image_20250727074723_ppp0.jpg"
$ws.Range("E11").Value = "PLACA_20250723145134"
$ws.Range("F11").Value = "Moura"
$ws.Range("G11").Value = 38.06587
$ws.Range("H11").Value = -7.221796
$ws.Range("I11").Value = "1490,161,1563,258"
$ws.Range("J11").Value = "'0.62"

# Row 18
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I18").Value = "1182,405,1231,455"
$ws.Range("J18").Value = "'0.76"
